{"js": "// Append 3 empty paragraphs followed by 5 text paragraphs to the end of the\n// document body (after the existing last paragraph, before the sectPr).\nconst body = context.document.body;\n\nconst newParagraphsText = [\n  \"\",\n  \"\",\n  \"\",\n  \"Fwrkthiuerwhtiurehihew\",\n  \"Rghwethyotueyuowuy\",\n  \"\\\\lerhguqkterytreh tljheroitqhrejklht;o \",\n  \"Phqrkehtu hrjtgkrh qotiorhtlr\",\n  \"Rtqrjhti oer jtoe j\",\n];\n\nfor (const text of newParagraphsText) {\n  body.insertParagraph(text, \"End\");\n}\n\nawait context.sync();\n", "ps1": "# Append 3 empty paragraphs followed by 5 text paragraphs to the end of the\n# document body (after the existing last paragraph, before the sectPr).\n$d = $word.ActiveDocument\n\nfunction Add-TrailingParagraph([string]$text) {\n    $d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null\n    if ($text -ne $null -and $text -ne \"\") {\n        $d.Paragraphs.Last.Range.Text = $text\n    }\n}\n\nAdd-TrailingParagraph \"\"\nAdd-TrailingParagraph \"\"\nAdd-TrailingParagraph \"\"\nAdd-TrailingParagraph \"Fwrkthiuerwhtiurehihew\"\nAdd-TrailingParagraph \"Rghwethyotueyuowuy\"\nAdd-TrailingParagraph \"\\lerhguqkterytreh tljheroitqhrejklht;o \"\nAdd-TrailingParagraph \"Phqrkehtu hrjtgkrh qotiorhtlr\"\nAdd-TrailingParagraph \"Rtqrjhti oer jtoe j\"\n"}
